# Fix the "card" column (A) on sheet "Card23" where rows were truncated to "2"
# instead of the correct value "23".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

$rows = @(3, 4, 5, 6, 7, 9, 10, 11, 12)
foreach ($r in $rows) {
    # Leading apostrophe forces the numeric-looking string to stay text,
    # matching the original inlineStr ("23" as text, not the number 23).
    $ws.Range("A$r").Value = "'23"
}
